$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap C4/D4 values
$ws.Range("C4").Value = "-"
$ws.Range("D4").Value = "Circuitos Elétricos 2"

# Swap C6/D6 values
$ws.Range("C6").Value = "-"
$ws.Range("D6").Value = "Circuitos Elétricos 2"
